# Auto-generated Excel COM-interop edit script.
# Scheduled-runner refresh of Leve currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) sourced from the latest market-board snapshot, across the
# ALC / ARM / BSM / CRP / CUL / GSM / LTW / WVR "Rafflesia" profit sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 10259.25
$ws.Range("I51").Value = 10666.667
$ws.Range("J51").Value = 10123.444
$ws.Range("K51").Value = 10666.667
$ws.Range("L51").Value = 10123.444
$ws.Range("M51").Value = -10182.667
$ws.Range("N51").Value = -11091.444
$ws.Range("H64").Value = 4495
$ws.Range("I64").Value = 4495
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 4495
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -4247
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 4495
$ws.Range("I67").Value = 4495
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 4495
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -3637
$ws.Range("N67").ClearContents()
$ws.Range("H132").Value = 38327.43
$ws.Range("I132").Value = 38327.43
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 114982.29
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -112452.29
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H61").Value = 3215.5
$ws.Range("J61").Value = 3444
$ws.Range("L61").Value = 3444
$ws.Range("N61").Value = -3868
$ws.Range("H136").Value = 3215.5
$ws.Range("J136").Value = 3444
$ws.Range("L136").Value = 10332
$ws.Range("N136").Value = -15432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1946.5454
$ws.Range("I20").Value = 1374.7142
$ws.Range("J20").Value = 2947.25
$ws.Range("K20").Value = 1374.7142
$ws.Range("L20").Value = 2947.25
$ws.Range("M20").Value = -1127.7142
$ws.Range("N20").Value = -3441.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 496.25
$ws.Range("I22").Value = 485
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 485
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -135
$ws.Range("N22").Value = -1200
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H86").Value = 5423
$ws.Range("I86").Value = 5402.3335
$ws.Range("J86").Value = 5454
$ws.Range("K86").Value = 5402.3335
$ws.Range("L86").Value = 5454
$ws.Range("M86").Value = -4279.3335
$ws.Range("N86").Value = -7700
$ws.Range("H89").Value = 5423
$ws.Range("I89").Value = 5402.3335
$ws.Range("J89").Value = 5454
$ws.Range("K89").Value = 27011.6675
$ws.Range("L89").Value = 27270
$ws.Range("M89").Value = -21395.6675
$ws.Range("N89").Value = -38502
$ws.Range("H122").Value = 1058.3334
$ws.Range("I122").Value = 1058.3334
$ws.Range("K122").Value = 3175.0002
$ws.Range("M122").Value = -725.0001999999999
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 10.166667
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 14.75
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 44.25
$ws.Range("M15").Value = 137
$ws.Range("N15").Value = -324.25
$ws.Range("H29").Value = 9.75
$ws.Range("J29").Value = 4
$ws.Range("L29").Value = 12
$ws.Range("N29").Value = -566
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 9000
$ws.Range("N80").Value = -10872
$ws.Range("M80").ClearContents()
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 27000
$ws.Range("N83").Value = -36360
$ws.Range("M83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1868.5
$ws.Range("I132").Value = 1555.75
$ws.Range("J132").Value = 2494
$ws.Range("K132").Value = 4667.25
$ws.Range("L132").Value = 7482
$ws.Range("M132").Value = -2137.25
$ws.Range("N132").Value = -12542

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 12800
$ws.Range("J14").Value = 12800
$ws.Range("L14").Value = 12800
$ws.Range("N14").Value = -13144
$ws.Range("H40").Value = 4000
$ws.Range("I40").Value = 4000
$ws.Range("K40").Value = 4000
$ws.Range("M40").Value = -3864
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 65000
$ws.Range("J75").Value = 65000
$ws.Range("L75").Value = 65000
$ws.Range("N75").Value = -66872
$ws.Range("H78").Value = 65000
$ws.Range("J78").Value = 65000
$ws.Range("L78").Value = 195000
$ws.Range("N78").Value = -204360
$ws.Range("H122").Value = 1597.5
$ws.Range("I122").Value = 1597.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4792.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2342.5
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 3800
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 11400
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -8930
$ws.Range("N126").Value = -19940
$ws.Range("H137").Value = 100000
$ws.Range("J137").Value = 100000
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200
